$wb = $excel.ActiveWorkbook

# Sheet "emp"
$ws = $wb.Worksheets.Item("emp")
$ws.Range("C2").Value = [double]"1.0143979550127285e-009"
$ws.Range("C3").Value = [double]"0.094049108711074747"
$ws.Range("C4").Value = [double]"0.010859487766057534"
$ws.Range("C5").Value = [double]"0.62702728760089754"

# Sheet "log_wage"
$ws = $wb.Worksheets.Item("log_wage")
$ws.Range("C2").Value = [double]"0.18534836829120924"
$ws.Range("C3").Value = [double]"0.65840609944649109"
$ws.Range("C4").Value = [double]"0.69861198658088808"
$ws.Range("C5").Value = [double]"0.70200611133569302"

# Sheet "log_hours"
$ws = $wb.Worksheets.Item("log_hours")
$ws.Range("C2").Value = [double]"0.042232394403873698"
$ws.Range("C3").Value = [double]"0.24541337938683686"
$ws.Range("C4").Value = [double]"0.80961077489224331"
$ws.Range("C5").Value = [double]"0.25996452292074584"
